# Apply the "Updated cryptos list" GitHub Actions refresh:
#  - Price (D) / Volume(1h) (E) columns get refreshed quote data for most rows.
#  - Rows 44-46 are re-ranked: OKB moves up to 44, ONDO to 45, Stacks to 46
#    (values for B/C/D/E on those rows follow the coin, not the row).
# All cells in this sheet are plain text (the source stores even numeric-looking
# prices like "579.11" or "1.00" as strings), so for any value that COM would
# otherwise auto-coerce to a number we force the cell to Text format first to
# keep it as a literal string (e.g. "1.00" must stay "1.00", not become 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "63.136.38"
$ws.Range("E2").Value = "  +2.45%  "
$ws.Range("D3").Value = "3.466.40"
$ws.Range("E3").Value = "  +2.16%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "579.09"
$ws.Range("D6").Value = "147.76"
$ws.Range("E6").Value = "  +3.18%  "
$ws.Range("D7").Value = "3.467.27"
$ws.Range("E7").Value = "  +2.20%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.479"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("E10").Value = "  +0.96%  "
$ws.Range("E11").Value = "  +1.93%  "
$ws.Range("D12").Value = "0.405"
$ws.Range("E12").Value = "  +5.27%  "
$ws.Range("D13").Value = "4.060.04"
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("D14").Value = "29.69"
$ws.Range("E14").Value = "  +6.21%  "
$ws.Range("E15").Value = "  +2.77%  "
$ws.Range("D16").Value = "3.467.63"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").Value = "63.154.70"
$ws.Range("E18").Value = "  +2.43%  "
$ws.Range("E19").Value = "  +3.56%  "
$ws.Range("E20").Value = "  +5.28%  "
$ws.Range("D21").Value = "9.28"
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("D22").Value = "389.08"
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("E23").Value = "  +2.12%  "
$ws.Range("D24").Value = "74.78"
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").Value = "3.610.61"
$ws.Range("E26").Value = "  +2.26%  "
$ws.Range("E27").Value = "  +1.93%  "
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("E29").Value = "  +3.14%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "8.17"
$ws.Range("E31").Value = "  +2.50%  "
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").Value = "1.38"
$ws.Range("E34").Value = "  -0.32%  "
$ws.Range("D35").Value = "23.64"
$ws.Range("E35").Value = "  +1.26%  "
$ws.Range("D36").Value = "5.32"
$ws.Range("E36").Value = "  +3.77%  "
$ws.Range("D37").Value = "7.08"
$ws.Range("E37").Value = "  +2.31%  "
$ws.Range("D38").Value = "32.04"
$ws.Range("E38").Value = "  +16.46%  "
$ws.Range("D39").Value = "170.16"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("E40").Value = "  +6.32%  "
$ws.Range("D41").Value = "3.502.40"
$ws.Range("E41").Value = "  +2.27%  "
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("D43").Value = "0.798"
$ws.Range("E43").Value = "  +2.14%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "42.32"
$ws.Range("E44").Value = "  -0.28%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").Value = "1.22"
$ws.Range("E45").Value = "  +4.73%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "1.73"
$ws.Range("E46").Value = "  +3.76%  "
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").Value = "2.628.57"
$ws.Range("E48").Value = "  +6.15%  "
$ws.Range("D49").Value = "2.28"
$ws.Range("E49").Value = "  +12.88%  "
$ws.Range("D50").Value = "23.03"
$ws.Range("E50").Value = "  +1.34%  "
$ws.Range("E51").Value = "  +2.37%  "
